# Daily attendance processing - 2026-01-23 12:00:54
#
# For a specific set of rows in the "Recorded By" column (G), reverse the
# order of the comma-separated list of recorders (e.g.
# "backup@backdoor.com, System" -> "System, backup@backdoor.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,8,28,29,30,31,32,33,34,54,55,56,57,58,59,60,80,81,82,87,106,107,108,113,132,133,134,139)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $current = $cell.Value2
    if ($current -ne $null) {
        $parts = $current -split ',\s*'
        $reversed = @()
        for ($i = $parts.Length - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $cell.Value = [string]::Join(", ", $reversed)
    }
}
